$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays formatted as text so values like
# "2.370.58" or "130.85" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Simple value updates (price/volume columns) ---
$ws.Range("D2").Value = "56.020.75"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.355.35"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "506.56"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "130.85"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").Value = "2.370.58"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "0.0975"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "4.82"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").Value = "0.322"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "2.773.40"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "55.967.60"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "21.54"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.385.12"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "9.96"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").Value = "312.16"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "4.03"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "65.38"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "0.371"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "0.146"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "171.71"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("E35").Value = "  -3.84%  "
$ws.Range("D36").Value = "17.71"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").Value = "3.66"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").Value = "36.17"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "125.67"
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").Value = "0.559"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").Value = "0.0896"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").Value = "242.59"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "0.0479"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "16.81"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "0.0207"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "16.75"
$ws.Range("E51").Value = "  -2.33%  "

# --- Row 30/31 swap: PancakeSwap <-> PEPE (with updated values) ---
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("E30").Value = "  -2.04%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.65"
$ws.Range("E31").Value = "  -0.06%  "
